# Edit script: rewrite SequenceConfig (sheet2) rows 2-49 with the new
# sequence data, and trim trailing-space command names on a handful of
# CommandList (sheet3) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SequenceConfig")

# Columns A and D frequently hold numeric-looking text ("5", "207", "750", ...).
# Force text format first so Excel keeps them as text (matching the
# "numberStoredAsText" cells in the source file) instead of silently
# converting them to real numbers.
$ws.Range("A2:A49").NumberFormat = "@"
$ws.Range("D2:D49").NumberFormat = "@"

$rows = @(
    @('5','ARMxl','stop_charge_session.sh','-'),
    @('5','AC_Source','set_output','OFF'),
    @('5','DC_Source','set_output','OFF'),
    @('5','AC_Source','set_voltage','207'),
    @('5','CLIM_chamber','write_setpoint','Temp 20'),
    @('5','CLIM_chamber','start_temp','-'),
    @('5','DC_Source','set_voltage','300'),
    @('5','DC_Source','set_c_limit','-160 6'),
    @('10','DC_Source','set_output','ON'),
    @('10','DC_Source','set_voltage ','750'),
    @('5','AC_Source','set_output','ON'),
    @('5','ARMxl','set_voltage_and_power.sh','8200 50'),
    @('5','ARMxl','start_charge_session.sh','-'),
    @('5','ARMxl','set_power.sh','940'),
    @('5400','CLIM_chamber','write_setpoint','Temp 25'),
    @('2760','CLIM_chamber','write_setpoint','Temp 35'),
    @('2760','CLIM_chamber','write_setpoint','Temp 50'),
    @('2760','CLIM_chamber','write_setpoint','Temp 55'),
    @('2760','CLIM_chamber','write_setpoint','Temp 60'),
    @('2760','CLIM_chamber','write_setpoint','Temp 65'),
    @('2760','CLIM_chamber','write_setpoint','Temp 70'),
    @('2760','CLIM_chamber','write_setpoint','Temp 75'),
    @('5','ARMxl','stop_charge_session.sh','-'),
    @('5','AC_Source','set_output','OFF'),
    @('5','DC_Source','set_output ','OFF'),
    @('5','AC_Source','usa_grid','-'),
    @('5','CLIM_chamber','write_setpoint','Temp 20'),
    @('5','CLIM_chamber','start_temp','-'),
    @('5','DC_Source','set_voltage ','300'),
    @('5','DC_Source','set_c_limit ','-160 6'),
    @('10','DC_Source','set_output ','ON'),
    @('10','DC_Source','set_voltage ','850'),
    @('5','AC_Source','set_output','ON'),
    @('5','ARMxl','set_voltage_and_power.sh','9200 50'),
    @('5','ARMxl','start_charge_session.sh','-'),
    @('5','ARMxl','set_power.sh','1050'),
    @('5400','CLIM_chamber','write_setpoint','Temp 25'),
    @('2760','CLIM_chamber','write_setpoint','Temp 35'),
    @('2760','CLIM_chamber','write_setpoint','Temp 50'),
    @('2760','CLIM_chamber','write_setpoint','Temp 55'),
    @('2760','CLIM_chamber','write_setpoint','Temp 60'),
    @('2760','CLIM_chamber','write_setpoint','Temp 65'),
    @('2760','CLIM_chamber','write_setpoint','Temp 70'),
    @('2760','CLIM_chamber','write_setpoint','Temp 75'),
    @('5','ARMxl','stop_charge_session.sh','-'),
    @('5','AC_Source','set_output','OFF'),
    @('5','DC_Source','set_output ','OFF'),
    @('5','CLIM_chamber','write_setpoint','Temp 25')
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# CommandList: strip stray trailing spaces from a few command-name cells.
$cmd = $wb.Worksheets.Item("CommandList")
$cmd.Cells.Item(15, 3).Value = "set_output"
$cmd.Cells.Item(15, 6).Value = "set_output"
$cmd.Cells.Item(16, 4).Value = "set_frequency"
$cmd.Cells.Item(16, 6).Value = "set_current"
$cmd.Cells.Item(17, 3).Value = "set_voltage"
$cmd.Cells.Item(17, 4).Value = "set_voltage"
$cmd.Cells.Item(18, 3).Value = "set_c_limit"
